# Remove the blank/index column (A), shifting columns B->A and C->B left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete()
